# Week 15 simulations added.
# Appends newly-simulated per-play numbers to the long simulation strings on
# the YDS and ST sheets, and updates the resulting aggregate stat cells on
# OFF, DEF, ST, TURNS and PEN to reflect the additional week of data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet: append new simulation numbers to the OFF/DEF R and P rows
# ---------------------------------------------------------------------
$ydsSheet = $wb.Worksheets.Item("YDS")

$ydsSheet.Range("B2").Value = $ydsSheet.Range("B2").Value2 + " 14 -3 10 22 1 29 18 4 11 5 18 23 2 0 7 8 3"
$ydsSheet.Range("C2").Value = $ydsSheet.Range("C2").Value2 + " 0 3 47 4 2 12 2 3 1 5 -1 8 2 1 1 2 9 1 -2 2 0 7 21 1 4 3 1"
$ydsSheet.Range("B3").Value = $ydsSheet.Range("B3").Value2 + " 1 10 11 5 3 25 24 21 9 5 6 10 12 6 1 13 9 -2 10 2 4 10 15 9 12 12 3 4 4 4 2 15 1 24 3"
$ydsSheet.Range("C3").Value = $ydsSheet.Range("C3").Value2 + " 6 15 5 7 5 9 20 8 16 23 18 5 13 16 13 7 5 20 5 6 24 4 18 0 8 5 11 10 6 -7 14 58"

# ---------------------------------------------------------------------
# OFF sheet: Home (row 2) / Road (row 3) totals
# ---------------------------------------------------------------------
$offSheet = $wb.Worksheets.Item("OFF")

$offSheet.Range("C2").Value = 119
$offSheet.Range("D2").Value = 8
$offSheet.Range("E2").Value = 8
$offSheet.Range("F2").Value = 57
$offSheet.Range("G2").Value = 40
$offSheet.Range("J2").Value = 26
$offSheet.Range("N2").Value = 15
$offSheet.Range("O2").Value = 15
$offSheet.Range("P2").Value = 6

$offSheet.Range("B3").Value = 7
$offSheet.Range("C3").Value = 194
$offSheet.Range("D3").Value = 3
$offSheet.Range("E3").Value = 24
$offSheet.Range("F3").Value = 90
$offSheet.Range("G3").Value = 25
$offSheet.Range("H3").Value = 27
$offSheet.Range("I3").Value = 33
$offSheet.Range("J3").Value = 52
$offSheet.Range("L3").Value = 266
$offSheet.Range("M3").Value = 180
$offSheet.Range("Q3").Value = 442

# ---------------------------------------------------------------------
# DEF sheet: Home (row 2) / Road (row 3) totals
# ---------------------------------------------------------------------
$defSheet = $wb.Worksheets.Item("DEF")

$defSheet.Range("C2").Value = 150
$defSheet.Range("D2").Value = 6
$defSheet.Range("F2").Value = 44
$defSheet.Range("G2").Value = 44
$defSheet.Range("J2").Value = 25
$defSheet.Range("N2").Value = 16
$defSheet.Range("O2").Value = 20

$defSheet.Range("C3").Value = 120
$defSheet.Range("D3").Value = 9
$defSheet.Range("E3").Value = 39
$defSheet.Range("F3").Value = 70
$defSheet.Range("G3").Value = 32
$defSheet.Range("H3").Value = 31
$defSheet.Range("I3").Value = 43
$defSheet.Range("J3").Value = 39
$defSheet.Range("L3").Value = 276
$defSheet.Range("M3").Value = 165
$defSheet.Range("Q3").Value = 447

# ---------------------------------------------------------------------
# ST sheet: aggregate cells + the four simulation strings
# ---------------------------------------------------------------------
$stSheet = $wb.Worksheets.Item("ST")

$stSheet.Range("B2").Value = 75
$stSheet.Range("D2").Value = 39
$stSheet.Range("F2").Value = 97
$stSheet.Range("G2").Value = 95
$stSheet.Range("J2").Value = 38
$stSheet.Range("K2").Value = 34

$stSheet.Range("B3").Value = 41

$stSheet.Range("D3").Value = $stSheet.Range("D3").Value2 + " 44 39 34 61 45 63"
$stSheet.Range("B4").Value = $stSheet.Range("B4").Value2 + " 56 61 56"
$stSheet.Range("D4").Value = $stSheet.Range("D4").Value2 + " 10 0 0 0 0 0"
$stSheet.Range("B5").Value = $stSheet.Range("B5").Value2 + " 15 33 7"
$stSheet.Range("D5").Value = $stSheet.Range("D5").Value2 + " 4 0 0 0 5"
$stSheet.Range("B6").Value = $stSheet.Range("B6").Value2 + " 29"

# ---------------------------------------------------------------------
# TURNS sheet: Road row (row 3)
# ---------------------------------------------------------------------
$turnsSheet = $wb.Worksheets.Item("TURNS")

$turnsSheet.Range("B3").Value = 9
$turnsSheet.Range("D3").Value = 9
$turnsSheet.Range("E3").Value = 10

# ---------------------------------------------------------------------
# PEN sheet: penalty counts/yards
# ---------------------------------------------------------------------
$penSheet = $wb.Worksheets.Item("PEN")

$penSheet.Range("D2").Value = 6
$penSheet.Range("B3").Value = 20
$penSheet.Range("D3").Value = 10
$penSheet.Range("D4").Value = 6
